# Add 8th test case scenario: Todo_Add_Issueses_Random_SelectDone
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24: new test case header (TC_NO = 8), matches the formatting used by
# rows 2/3 (A: centered number, C/D: wrap-text description cells)
# -4108 = xlCenter (used for both HorizontalAlignment and VerticalAlignment)
$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 1).HorizontalAlignment = -4108
$ws.Cells.Item(24, 1).VerticalAlignment = -4108

$ws.Cells.Item(24, 2).Value = "Todo_Add_Issueses_Random_SelectDone"

$ws.Cells.Item(24, 3).Value = "* https://todomvc.com/examples/vue/#   adresine girilir"
$ws.Cells.Item(24, 3).WrapText = $true

$ws.Cells.Item(24, 4).Value = "Sayfanın Başarılı şekilde açıldıgı görülür"
$ws.Cells.Item(24, 4).WrapText = $true

# Row 25: steps / expected result
$ws.Cells.Item(25, 3).Value = "* What needs to be done ? İnput'u içerisine herhangi bir değer girilip Enter'a basılır bu islem birden fazla olucak sekilde defa tekrarlanir"
$ws.Cells.Item(25, 3).WrapText = $true

$ws.Cells.Item(25, 4).Value = "Yazilan  değerler başarılı şekilde eklendiği görülür,  X item left yazisi görülür , All Active Completed butonları görülür || Eklenen deger active olarak eklendigi gorulur "
$ws.Cells.Item(25, 4).WrapText = $true
$ws.Rows.Item(25).RowHeight = 45

# Row 26: additional step / expected result
$ws.Cells.Item(26, 3).Value = "* Rastgele 1 issue'nun yanındaki checkBox tiklanir"
$ws.Cells.Item(26, 3).WrapText = $true

$ws.Cells.Item(26, 4).Value = "Tıklanilan issue DONE statüsüne geldigi gorulur"
$ws.Cells.Item(26, 4).WrapText = $true

# Update selection / view to match the authored state
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 16
